# Commit: "add fr jurisdiction, snomed parameters, ..."
#
# The ValueSet's "Metadata" sheet lists Property/Value pairs. This edit:
#   - refreshes the generation "Date" value (B8)
#   - fills in the previously-blank "Jurisdiction" value with "FRANCE" (B11)

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

$metadata.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$metadata.Range("B11").Value = "FRANCE"
